$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.606.55"
$ws.Range("E2").Value = "  +0.18%  "

# Row 3
$ws.Range("D3").Value = "1.962.31"
$ws.Range("E3").Value = "  +0.33%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "244.86"
$ws.Range("E5").Value = "  +0.13%  "

# Row 6
$ws.Range("E6").Value = "  -0.66%  "

# Row 7
$ws.Range("D7").Value = "58.70"
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").Value = "0.378"
$ws.Range("E9").Value = "  +2.64%  "

# Row 10
$ws.Range("D10").Value = "0.0806"
$ws.Range("E10").Value = "  -6.47%  "

# Row 11
$ws.Range("E11").Value = "  -0.79%  "

# Row 12
$ws.Range("E12").Value = "  -1.13%  "

# Row 13
$ws.Range("D13").Value = "0.832"

# Row 14
$ws.Range("D14").Value = "2.249.45"
$ws.Range("E14").Value = "  +0.32%  "

# Row 15
$ws.Range("D15").Value = "13.76"
$ws.Range("E15").Value = "  +0.32%  "

# Row 16
$ws.Range("E16").Value = "  +0.92%  "

# Row 17
$ws.Range("D17").Value = "1.955.44"
$ws.Range("E17").Value = "  -0.73%  "

# Row 18
$ws.Range("D18").Value = "36.535.41"
$ws.Range("E18").Value = "  +0.16%  "

# Row 19
$ws.Range("D19").Value = "69.84"
$ws.Range("E19").Value = "  -0.42%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0856"
$ws.Range("E20").Value = "  -2.85%  "

# Row 21
$ws.Range("D21").Value = "228.88"
$ws.Range("E21").Value = "  -0.56%  "

# Row 22
$ws.Range("D22").Value = "5.06"
$ws.Range("E22").Value = "  -0.75%  "

# Row 23
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.23%  "

# Row 24
$ws.Range("D24").Value = "2.46"
$ws.Range("E24").Value = "  -1.19%  "

# Row 25
$ws.Range("E25").Value = "  +1.54%  "

# Row 26
$ws.Range("D26").Value = "9.26"
$ws.Range("E26").Value = "  -1.89%  "

# Row 27
$ws.Range("E27").Value = "  +1.40%  "

# Row 28
$ws.Range("D28").Value = "160.48"
$ws.Range("E28").Value = "  -1.29%  "

# Row 29
$ws.Range("D29").Value = "19.45"
$ws.Range("E29").Value = "  -1.08%  "

# Row 30
$ws.Range("E30").Value = "  +1.22%  "

# Row 31
$ws.Range("E31").Value = "  -3.33%  "

# Row 33
$ws.Range("D33").Value = "0.0621"
$ws.Range("E33").Value = "  -3.51%  "

# Row 34
$ws.Range("D34").Value = "4.33"
$ws.Range("E34").Value = "  +0.21%  "

# Row 35
$ws.Range("E35").Value = "  -0.14%  "

# Row 36
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "2.24"
$ws.Range("E36").Value = "  +2.31%  "

# Row 37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "3.39"
$ws.Range("E37").Value = "  +11.08%  "

# Row 39
$ws.Range("D39").Value = "5.77"
$ws.Range("E39").Value = "  -10.25%  "

# Row 40
$ws.Range("D40").Value = "0.0982"
$ws.Range("E40").Value = "  -1.93%  "

# Row 41
$ws.Range("D41").Value = "2.91"
$ws.Range("E41").Value = "  +0.60%  "

# Row 42
$ws.Range("E42").Value = "  -1.52%  "

# Row 43
$ws.Range("E43").Value = "  +0.00%  "

# Row 44
$ws.Range("D44").Value = "16.03"
$ws.Range("E44").Value = "  -1.16%  "

# Row 45
$ws.Range("D45").Value = "1.366.03"
$ws.Range("E45").Value = "  +0.62%  "

# Row 46
$ws.Range("E46").Value = "  -0.92%  "

# Row 47
$ws.Range("D47").Value = "87.99"
$ws.Range("E47").Value = "  -0.77%  "

# Row 48
$ws.Range("E48").Value = "  -1.25%  "

# Row 49
$ws.Range("D49").Value = "2.83"
$ws.Range("E49").Value = "  +0.33%  "

# Row 50
$ws.Range("D50").Value = "2.139.99"
$ws.Range("E50").Value = "  +0.23%  "

# Row 51
$ws.Range("D51").Value = "43.79"
$ws.Range("E51").Value = "  -5.11%  "
